# Applies the ICDC TCL01 test edit: updates the Cypher queries/expected
# filenames on the existing CasesTab/SamplesTab/FilesTab rows and appends a
# new StudyFilesTab row (row 5) with matching styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Query / filename text blocks (single-quoted here-strings so backticks and
# `$` in the Cypher text are taken literally).
# ---------------------------------------------------------------------------

$statQueryV1 = @'
MATCH (p:program)<--(s:study)<--(c)
MATCH (cf)-->(samp:sample)
WHERE samp.specific_sample_pathology IN ['Lymphoma']
MATCH (cf:file)-[*]->(c:case)
OPTIONAL MATCH (sf:file)-->(s)
RETURN
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
   count(distinct samp) AS Samples,
    count(distinct cf) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
MATCH (samp:sample)-->(c)
WHERE  samp.specific_sample_pathology in ['Lymphoma']
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
  coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
order by c.case_id asc
limit 100
'@

$neo4jFile = 'TC03_Canine_Filter_SamplePatho-Lymphoma_Neo4jData.xlsx'
$webFile = 'TC03_Canine_Filter_SamplePatho-Lymphoma_WebData.xlsx'

$statQueryV2 = @'
MATCH (p:program)<--(s:study)<--(c)
MATCH (cf)-->(samp:sample)
WHERE samp.specific_sample_pathology IN ['Lymphoma']
MATCH (cf:file)-[*]->(c:case)
OPTIONAL MATCH (sf:file)-->(s)
RETURN
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct cf) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$neo4jFileTypo = 'TC03_Canine_Filter_SamplePatho-Lymohoma_Neo4jData.xlsx'

$studyFilesQuery = @'
MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(diag:diagnosis)
MATCH (c)<--(demo:demographic)
MATCH (samp:sample)-->(c)
WHERE samp.specific_sample_pathology IN ['Lymphoma']
WITH
        DISTINCT f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
  order by 'File Name' asc
  limit 100
'@

$samplesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE  samp.specific_sample_pathology IN ['Lymphoma']
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed,
        coalesce(diag.disease_term,'') AS Diagnosis, 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
order by samp.sample_id asc
limit 100
'@

$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
MATCH (f)-->(samp:sample)
WHERE samp.specific_sample_pathology IN ['Lymphoma']
 MATCH (f)-[*]->(samp:sample)
WITH
        DISTINCT f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN
        coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_type, '') AS `File Type`,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
        order by f.file_name asc
        limit 100
'@

# ---------------------------------------------------------------------------
# Row 2 (CasesTab)
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = $casesQuery
$ws.Range("B2").WrapText = $true

$ws.Range("C2").Value = $statQueryV1
$ws.Range("C2").WrapText = $true
$ws.Range("C2").VerticalAlignment = -4108

$ws.Range("D2").Value = $neo4jFile
$ws.Range("E2").Value = $webFile

$ws.Rows.Item(2).RowHeight = 304.5

# ---------------------------------------------------------------------------
# Row 3 (SamplesTab)
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = $samplesQuery
$ws.Range("B3").WrapText = $true
$ws.Range("B3").Font.Size = 18

$ws.Range("C3").Value = $statQueryV2
$ws.Range("C3").WrapText = $true
$ws.Range("C3").VerticalAlignment = -4108

$ws.Range("D3").Value = $neo4jFileTypo
$ws.Range("E3").Value = $webFile

$ws.Rows.Item(3).RowHeight = 409.5

# ---------------------------------------------------------------------------
# Row 4 (FilesTab)
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = $filesQuery
$ws.Range("B4").WrapText = $true
$ws.Range("B4").Font.Size = 18

$ws.Range("C4").Value = $statQueryV2
$ws.Range("C4").WrapText = $true
$ws.Range("C4").VerticalAlignment = -4108

$ws.Range("D4").Value = $neo4jFile
$ws.Range("E4").Value = $webFile

$ws.Rows.Item(4).RowHeight = 409.5

# ---------------------------------------------------------------------------
# Row 5 (new StudyFilesTab row)
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "StudyFilesTab"

$ws.Range("B5").Value = $studyFilesQuery
$ws.Range("B5").WrapText = $true
$ws.Range("B5").Font.Size = 18

$ws.Range("C5").Value = $statQueryV2
$ws.Range("C5").WrapText = $true
$ws.Range("C5").Font.Size = 18

$ws.Range("D5").Value = $neo4jFile
$ws.Range("E5").Value = $webFile

$ws.Rows.Item(5).RowHeight = 409.5

# ---------------------------------------------------------------------------
# View state: selection on C5 (matches the author's last edit position)
# ---------------------------------------------------------------------------
$ws.Range("C5").Select()
